# Generate Report for Handback
#
# This script brings the "localization-status" report up to date after a
# handback event: the md files' status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the newly-produced handback target
# files are recorded (Latest Target File / Latest Handback File columns),
# their datetime stamps are updated, and the columns that now hold longer
# file names are widened.

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30c1cd47a6d2faab3779d9d65d302346a2e1376d/e2e/9c60a244-244f-44a5-a2b5-ab4e50383cae.md"
$mdName1 = "9c60a244-244f-44a5-a2b5-ab4e50383cae.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30c1cd47a6d2faab3779d9d65d302346a2e1376d/e2e/e25db826-d4aa-4702-934f-31a83c67c0f4.md"
$mdName2 = "e25db826-d4aa-4702-934f-31a83c67c0f4.md"

function Update-LanguageSheet($ws, $targetFile2, $targetFile3, $handbackFile2, $handbackFile3, $handbackDateTime2, $handbackDateTime3) {
    # Status column (C) -> handed back
    if ($ws.Range("C2").Text -eq $statusOld) { $ws.Range("C2").Value = $statusNew }
    if ($ws.Range("C3").Text -eq $statusOld) { $ws.Range("C3").Value = $statusNew }

    # Latest Target File / Latest Handback File / Latest Handback DateTime
    $ws.Range("I2").Value = $targetFile2
    $ws.Range("J2").Value = $handbackFile2
    $ws.Range("K2").Value = $handbackDateTime2

    $ws.Range("I3").Value = $targetFile3
    $ws.Range("J3").Value = $handbackFile3
    $ws.Range("K3").Value = $handbackDateTime3

    # Recreate the hyperlinks so that the new "Latest Target File" cells
    # (I2 / I3) point back at the source markdown files, same as column A.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $mdName2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrl2, "", "", $mdName2)

    # Widen the columns that now contain the longer generated file names.
    $ws.Columns.Item(3).ColumnWidth = 29.17
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}

$wb = $excel.ActiveWorkbook

### ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
if ($wsOverview.Range("E2").Text -eq $statusOld) { $wsOverview.Range("E2").Value = $statusNew }
if ($wsOverview.Range("F2").Text -eq $statusOld) { $wsOverview.Range("F2").Value = $statusNew }
if ($wsOverview.Range("E3").Text -eq $statusOld) { $wsOverview.Range("E3").Value = $statusNew }
if ($wsOverview.Range("F3").Text -eq $statusOld) { $wsOverview.Range("F3").Value = $statusNew }

$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

### ---- per-language sheets (zh-cn / de-de) ------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $wsZhCn `
    "9c60a244-244f-44a5-a2b5-ab4e50383cae.md" `
    "e25db826-d4aa-4702-934f-31a83c67c0f4.md" `
    "9c60a244-244f-44a5-a2b5-ab4e50383cae.2272ef5e6c801a16ddf307674235daf01b26af76.zh-cn.xlf" `
    "e25db826-d4aa-4702-934f-31a83c67c0f4.e924ef4ee79a4e9af59009b08e7911d899d73e63.zh-cn.xlf" `
    "2016-08-18 13:02:31" `
    "2016-08-18 13:02:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $wsDeDe `
    "9c60a244-244f-44a5-a2b5-ab4e50383cae.md" `
    "e25db826-d4aa-4702-934f-31a83c67c0f4.md" `
    "9c60a244-244f-44a5-a2b5-ab4e50383cae.2272ef5e6c801a16ddf307674235daf01b26af76.de-de.xlf" `
    "e25db826-d4aa-4702-934f-31a83c67c0f4.e924ef4ee79a4e9af59009b08e7911d899d73e63.de-de.xlf" `
    "2016-08-18 13:02:39" `
    "2016-08-18 13:02:39"
